$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows 2-44: Price (D) and Volume(1h) (E) changes ---
$ws.Range("D2").Value = "'26.543.36"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +2.25%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'1.682.31"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  +2.74%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  -0.19%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'217.36"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Value = "'0.5337"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  +3.52%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("E7").Value = "'  -0.22%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'0.2677"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  +4.77%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.06421"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  +3.17%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'21.48"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  +6.10%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.07790"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  +3.36%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'4.509"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  +3.66%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'1.674.63"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  +2.19%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'0.5630"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  +4.44%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'0.0₅8432"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  +6.66%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("E16").Value = "'  +2.13%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'26.583.05"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  +2.35%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("E18").Value = "'  -0.21%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'4.798"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  +3.67%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'195.73"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  +6.03%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("E21").Value = "'  +4.62%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'6.374"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  +4.95%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'1.001"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  -0.25%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'143.31"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  -1.31%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'0.1283"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  +8.18%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'7.475"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  +2.34%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("E27").Value = "'  +4.81%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'1.419"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  +3.36%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'0.06124"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  +3.15%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'1.279"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  +2.96%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'3.607"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  +7.88%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("E32").Value = "'  +3.63%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'1.707"
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").Value = "'1.016"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  +5.16%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'2.417"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  +1.52%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'2.786"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  +1.93%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'0.5710"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  -1.97%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'0.01643"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  +3.19%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'5.958"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  +4.83%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'0.8739"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  +3.99%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'1.061.63"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  +2.25%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'1.001"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  -0.11%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'100.03"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  +0.50%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'1.833.27"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  +2.32%  "
$ws.Range("E44").Style = "Normal"

# --- Row 45: new coin BabyDogeCoin inserted, rows 45-50 (old) shift to 46-51, Algorand (old 51) drops off ---
$ws.Range("B45").Value = "'BabyDogeCoin"
$ws.Range("B45").Style = "Normal"
$ws.Range("C45").Value = "'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("C45").Style = "Normal"
$ws.Range("D45").Value = "'0.0₈113"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  +6.16%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("B46").Value = "'Aave"
$ws.Range("B46").Style = "Normal"
$ws.Range("C46").Value = "'https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("C46").Style = "Normal"
$ws.Range("D46").Value = "'57.26"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  +5.56%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("B47").Value = "'EnergySwap"
$ws.Range("B47").Style = "Normal"
$ws.Range("C47").Value = "'https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("C47").Style = "Normal"
$ws.Range("D47").Value = "'8.162"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  +3.10%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("B48").Value = "'Frax"
$ws.Range("B48").Style = "Normal"
$ws.Range("C48").Value = "'https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("C48").Style = "Normal"
$ws.Range("D48").Value = "'0.9989"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  -0.17%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("B49").Value = "'Cronos"
$ws.Range("B49").Style = "Normal"
$ws.Range("C49").Value = "'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("C49").Style = "Normal"
$ws.Range("D49").Value = "'0.05203"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  +0.31%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("B50").Value = "'Aptos"
$ws.Range("B50").Style = "Normal"
$ws.Range("C50").Value = "'https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("C50").Style = "Normal"
$ws.Range("D50").Value = "'6.091"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  +5.67%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("B51").Value = "'Mantle"
$ws.Range("B51").Style = "Normal"
$ws.Range("C51").Value = "'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("C51").Style = "Normal"
$ws.Range("D51").Value = "'0.4239"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  +0.15%  "
$ws.Range("E51").Style = "Normal"
